# Revert "Added Pyrite, 50 Bq FHY pH9"
#
# This undoes the earlier commit that:
#   - appended 15 "Pyrite" data rows (rows 62-76)
#   - added a "Pyrite" shared string
#   - overwrote row 19 (50 Bq Goethite) and row 21 (500 Bq Goethite) with
#     recomputed values, and added a comment to row 19
#
# We restore the pre-commit state: drop the Pyrite rows, put row 19 / row 21
# back to their original numbers, and mark row 19 with its original
# "Need solid counts" comment.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Remove the 15 appended "Pyrite" rows (A62:N76) ---------------------
$ws.Range("A62:N76").EntireRow.Delete()

# --- 2. Restore row 19 (50 Bq/mL, Goethite) values --------------------------
$ws.Cells.Item(19, 3).Value2  = 0                     # C19
$ws.Cells.Item(19, 4).Value2  = 0                     # D19
$ws.Cells.Item(19, 8).Value2  = 1                     # H19
$ws.Cells.Item(19, 9).Value2  = 0                     # I19
$ws.Cells.Item(19, 12).Value2 = 42481                 # L19 (Date of Input)
$ws.Cells.Item(19, 13).Value2 = $false                # M19 (Include?)
$ws.Cells.Item(19, 14).Value2 = "Need solid counts"   # N19 (Comments)

# --- 3. Restore row 21 (500 Bq/mL, Goethite) values --------------------------
# (written in plain decimal -- this PowerShell parser doesn't accept `E`
# scientific-notation numeric literals)
$ws.Cells.Item(21, 3).Value2  = 0.04408467452874436   # C21
$ws.Cells.Item(21, 4).Value2  = 0.029860179390728075  # D21
$ws.Cells.Item(21, 5).Value2  = 13998.532901237906    # E21
$ws.Cells.Item(21, 6).Value2  = 1891.3936416951028    # F21
$ws.Cells.Item(21, 8).Value2  = 0.8396335735958603    # H21
$ws.Cells.Item(21, 9).Value2  = 0.11344600278165595   # I21
$ws.Cells.Item(21, 10).Value2 = 8.935                 # J21
$ws.Cells.Item(21, 11).Value2 = 0.04725815626252589   # K21

# --- 4. Restore the view selection / scroll state ---------------------------
$ws.Range("N20").Select()
$win = $excel.ActiveWindow
$win.ScrollColumn = 9
$win.ScrollRow = 2
